$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Memory Map")

# Row 9: A9 was formula =2*4096 (8192) -> now plain value 2048
$ws.Range("A9").Value = 2048

# Row 11: A11 was plain value 2048 -> now plain value 8192
$ws.Range("A11").Value = 8192

# Swap the G/I descriptive text between row 9 and row 11
$ws.Range("G9").Value = "Local variables"
$ws.Range("I9").Value = "Access to subroutine and exception stacks (current task)"
$ws.Range("G11").Value = "Character RAM"
$ws.Range("I11").Value = "Access to the character RAM in the CPU address space "

$ws.Range("L9").Font.Size = 12
$ws.Range("L10").Font.Size = 12
$ws.Range("L11").Font.Size = 12

$ws.Range("J13").Select() | Out-Null

# Column G gets slightly narrower; column I becomes hidden
$ws.Columns.Item(7).ColumnWidth = 33.33
$ws.Columns.Item(9).Hidden = $true

# Add a new "Sheet2" after the last sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheet2 = $wb.Worksheets.Add($null, $lastSheet)
$sheet2.Name = "Sheet2"
$sheet2.Range("B2").Value = "03B000"
$sheet2.Range("C2").Formula = "=HEX2BIN(B2/2^11,7)"
$sheet2.Columns.Item(3).ColumnWidth = 12.2
$sheet2.Range("B2").Select() | Out-Null

# Restore "Memory Map" as the active sheet/tab
$ws.Activate()

Write-Host "done"
